# Trading update: 2026-02-17 07:54:54
# Append the newest trade (Trade #19, still OPEN) as row 20 on both the
# "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 20

    # Numeric columns - plain values are fine.
    $ws.Range("A$row").Value = 19
    $ws.Range("F$row").Value = 0.01
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 99.94513875800263
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("Q$row").Value = 0

    # Plain text columns - not date/time-shaped, safe to assign directly.
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"

    # Date/time-shaped text columns - a direct .Value assignment gets
    # auto-parsed into a date/time serial by Excel's smart input. Build the
    # literal text via a formula, then copy / paste-values so the cell ends
    # up holding plain text (matching the source data export), not a date.
    $ws.Range("B$row").Formula = "=""2026-02-17"""
    $ws.Range("B$row").Copy()
    $ws.Range("B$row").PasteSpecial(-4163)

    $ws.Range("C$row").Formula = "=""07:54:39"""
    $ws.Range("C$row").Copy()
    $ws.Range("C$row").PasteSpecial(-4163)

    # G (Exit Price) and P (Exit Reason) stay blank - the trade is still OPEN.
}

Write-Output "done"
